$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New player row (row 32): Santiago Sandoval -----------------------------
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "Santiago Sandoval"
$ws.Range("C32").Value = 18
$ws.Range("D32").Value = 56
$ws.Range("E32").Value = 165
$ws.Range("F32").Value = "Mediocampista"
$ws.Range("G32").Value = "Medio Ofensivo"

# Foto_URL (H32) carries a hyperlink, same as the rest of the column.
$ws.Range("H32").Value = "https://i.imgur.com/VwmqSSM.png"
$ws.Hyperlinks.Add($ws.Range("H32"), "https://i.imgur.com/VwmqSSM.png")
$ws.Range("H32").Style = "Hipervínculo"

# Carnet_URL (I32) is plain text, no hyperlink (matches existing sheet pattern).
$ws.Range("I32").Value = "https://cdn.resfu.com/img_data/players/medium/3348548.jpg?size=120x&lossy=1"

# --- View state: selection moved to H33 while scrolled further down ---------
$ws.Range("H33").Select()
